$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (RMI 2025 Summer Recovery note): update deployment note text with
# fuller names, and grow the row to fit the longer note. The cell already
# carries the wrap-text style (s=3), so a direct value write keeps the style.
$ws.Range("F6").Value = "Used CTD/rosette with no issues, People: Nicole Waite, Lori Garzio, Dave Aragon, MOO students, Seatow vessel Jennie Lee.`nNOTE on the DEP glider ru32 that was deployed with ru39 - ru32 had issues and was recovered by Sea Tow without RU staff so there are no water samples for that glider recovery"
$ws.Rows.Item(6).RowHeight = 64

# --- Row 7 (new): RMI 2025 Fall Deployment
# Clone formatting from row 2 (same shape: deployment row) so the date/note
# cells pick up the existing number-format / wrap-text styles instead of
# minting new ones.
$ws.Range("A2:F2").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122)
$ws.Range("A7").Value = "RMI"
$ws.Range("B7").Value = "RMI 2025 Fall Deployment"
$ws.Range("C7").Value = "ru39-20251024T1506"
$ws.Range("D7").Value = "deployment"
$ws.Range("E7").Value = 45954
$ws.Range("F7").Value = "Used single niskin to collect water samples at 2m and 8m. People: Brian Buckingham, Lori Garzio, Ashley Hann. Did the first cast, then went to remove tape from ru43 fluorometer, then continued with casts. Rough seas so we poisoned samples at the dock"
$ws.Rows.Item(7).RowHeight = 48

# --- Row 8 (new): RMI 2025 Fall Recovery
# Clone formatting from row 3 (same shape: recovery row).
$ws.Range("A3:F3").Copy()
$ws.Range("A8:F8").PasteSpecial(-4122)
$ws.Range("A8").Value = "RMI"
$ws.Range("B8").Value = "RMI 2025 Fall Recovery"
$ws.Range("C8").Value = "ru39-20251024T1506"
$ws.Range("D8").Value = "recovery"
$ws.Range("E8").Value = 45982
$ws.Range("F8").Value = "Used rosette to collect water samples at 2m and 8m. For the first cast, the bottle for the 2m sample didn't close due to a mechanical issue, so immediately put the rosette back in the water to collect the surface sample. We consider this one cast since they were so close together. Very calm seas. People: Nicole Waite, Lori Garzio, Ashley Hann on the SeaTow Jennie Lee"
$ws.Rows.Item(8).RowHeight = 64

# --- View state: move the frozen-pane scroll position down and select the
# new last data cell, matching where the author's cursor ended up.
$ws.Range("F8").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
